# Insert a new "test_suite" sheet before the first sheet (AddCustomerTest),
# populate it with a TCID / Runmode table, add a "runmode" column to the
# AddCustomerTest sheet, and make AddCustomerTest the active tab again.

$wb = $excel.ActiveWorkbook

# --- New "test_suite" sheet, inserted before the current first sheet ---
$firstSheet = $wb.Worksheets.Item(1)
$testSuite = $wb.Worksheets.Add($firstSheet)
$testSuite.Name = "test_suite"

$testSuite.Range("A1").Value = "TCID"
$testSuite.Range("B1").Value = "Runmode"
$testSuite.Range("A2").Value = "BankManagerLoginTest"
$testSuite.Range("A3").Value = "AddCustomerTest"
$testSuite.Range("A4").Value = "OpenAccountTest"
$testSuite.Range("B2").Value = "Y"
$testSuite.Range("B3").Value = "Y"
$testSuite.Range("B4").Value = "Y"
[void]$testSuite.Range("B4").Select()
$testSuite.PageSetup.Orientation = 1

# --- Add a "runmode" column (E) to AddCustomerTest ---
$addCustomer = $wb.Worksheets.Item("AddCustomerTest")
$addCustomer.Range("E1").Value = "runmode"
$addCustomer.Range("E2").Value = "Y"
$addCustomer.Range("E3").Value = "N"
$addCustomer.Range("E4").Value = "Y"
$addCustomer.Range("E5").Value = "Y"
[void]$addCustomer.Range("F19").Select()

# AddCustomerTest remains the active/visible tab
$addCustomer.Activate()
